# Update regression output table (hourly regression with fixed effects)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1064960414663828

$ws.Range("B3").Value = 0.001360557758454122
$ws.Range("C3").Value = 0.0006657543448392394
$ws.Range("D3").Value = 2.240979161002647
$ws.Range("E3").Value = 0.09708252463288844
$ws.Range("F3").Value = 0.00005569888555048592
$ws.Range("G3").Value = 0.002665416631357759
$ws.Range("H3").Value = 0.1078565992248369

$ws.Range("B4").Value = 0.01152617729072993
$ws.Range("C4").Value = 0.001405153811283492
$ws.Range("D4").Value = 8.071191060019499
$ws.Range("E4").Value = 0.07338213104006673
$ws.Range("F4").Value = 0.008772118104493524
$ws.Range("G4").Value = 0.01428023647696634
$ws.Range("H4").Value = 0.1180222187571127

$ws.Range("B5").Value = 0.01373047161822737
$ws.Range("C5").Value = 0.006491186634822179
$ws.Range("D5").Value = 6.633884970403383
$ws.Range("E5").Value = 0.1457911174149475
$ws.Range("F5").Value = 0.001007939167300922
$ws.Range("G5").Value = 0.02645300406915381
$ws.Range("H5").Value = 0.1202265130846101

$ws.Range("B6").Value = 0.01164646475784711
$ws.Range("C6").Value = 0.003765595632747634
$ws.Range("D6").Value = 5.178272378650155
$ws.Range("E6").Value = 0.07686960457822326
$ws.Range("F6").Value = 0.004266012890411725
$ws.Range("G6").Value = 0.01902691662528249
$ws.Range("H6").Value = 0.1181425062242299

$ws.Range("B7").Value = 0.01481699108592017
$ws.Range("C7").Value = 0.006483547670881736
$ws.Range("D7").Value = 4.008740890156023
$ws.Range("E7").Value = 0.1140347679492291
$ws.Range("F7").Value = 0.002109438194913143
$ws.Range("G7").Value = 0.0275245439769272
$ws.Range("H7").Value = 0.121313032552303

$ws.Range("B8").Value = 0.01198874003064032
$ws.Range("C8").Value = 0.006084024219851028
$ws.Range("D8").Value = 4.591842971072017
$ws.Range("E8").Value = 0.1205668699114027
$ws.Range("F8").Value = 0.00006424008068581881
$ws.Range("G8").Value = 0.02391323998059483
$ws.Range("H8").Value = 0.1184847814970231

$ws.Range("B9").Value = 0.009717951609461059
$ws.Range("C9").Value = 0.007964177051624503
$ws.Range("D9").Value = 4.320388793747422
$ws.Range("E9").Value = 0.1032628555774512
$ws.Range("F9").Value = -0.005891593498264538
$ws.Range("G9").Value = 0.02532749671718666
$ws.Range("H9").Value = 0.1162139930758438

$ws.Range("B10").Value = -0.1064960414663828
$ws.Range("C10").Value = 0.0004814384452962593
$ws.Range("D10").Value = -239.2625494043926
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.1074396465421868
$ws.Range("G10").Value = -0.1055524363905787

$ws.Range("B11").Value = -0.04758467142744626
$ws.Range("C11").Value = 0.0005251722608214245
$ws.Range("D11").Value = -97.05769047497839
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -0.04861399349439435
$ws.Range("G11").Value = -0.04655534936049816
$ws.Range("H11").Value = 0.05891137003893652

$ws.Range("B12").Value = -0.03600997783196849
$ws.Range("C12").Value = 0.0005069762511751827
$ws.Range("D12").Value = -77.57845744213114
$ws.Range("E12").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000003896347234450767
$ws.Range("F12").Value = -0.03700363625136938
$ws.Range("G12").Value = -0.03501631941256762
$ws.Range("H12").Value = 0.07048606363441429

$ws.Range("B13").Value = -0.03312148630382851
$ws.Range("C13").Value = 0.0005078854039619428
$ws.Range("D13").Value = -70.57631248675035
$ws.Range("E13").Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000008624017266256516
$ws.Range("F13").Value = -0.03411692664755753
$ws.Range("G13").Value = -0.03212604596009948
$ws.Range("H13").Value = 0.07337455516255427

$ws.Range("B14").Value = -0.02711818916936227
$ws.Range("C14").Value = 0.0004948760748773819
$ws.Range("D14").Value = -61.03429212966432
$ws.Range("E14").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000113241562069149
$ws.Range("F14").Value = -0.02808813160792204
$ws.Range("G14").Value = -0.02614824673080249
$ws.Range("H14").Value = 0.07937785229702052

$ws.Range("B15").Value = -0.02232895710583033
$ws.Range("C15").Value = 0.0004886076926807012
$ws.Range("D15").Value = -50.05908973490114
$ws.Range("E15").Value = 0.000000000000000000000000000000000000000000000000000000002793410367242538
$ws.Range("F15").Value = -0.02328661370273277
$ws.Range("G15").Value = -0.0213713005089279
$ws.Range("H15").Value = 0.08416708436055245

$ws.Range("B16").Value = -0.02154908315537955
$ws.Range("C16").Value = 0.0004774112571119029
$ws.Range("D16").Value = -48.58678148271015
$ws.Range("E16").Value = 0.000000000000000000000000000000000000002836905899859235
$ws.Range("F16").Value = -0.02248479507410671
$ws.Range("G16").Value = -0.02061337123665237
$ws.Range("H16").Value = 0.08494695831100324

$ws.Range("B17").Value = -0.01989460376376301
$ws.Range("C17").Value = 0.0004846323122433717
$ws.Range("D17").Value = -45.63276626024594
$ws.Range("E17").Value = 0.000000000000000000000000000000000000387258275323383
$ws.Range("F17").Value = -0.02084446873002404
$ws.Range("G17").Value = -0.01894473879750198
$ws.Range("H17").Value = 0.08660143770261977

$ws.Range("B18").Value = -0.01768384768455453
$ws.Range("C18").Value = 0.0004927594257718047
$ws.Range("D18").Value = -39.52624895587892
$ws.Range("E18").Value = 0.0000000000000000000000000007951326698258996
$ws.Range("F18").Value = -0.01864964154869635
$ws.Range("G18").Value = -0.01671805382041271
$ws.Range("H18").Value = 0.08881219378182825

$ws.Range("B19").Value = -0.01501933554117781
$ws.Range("C19").Value = 0.0004870935315581139
$ws.Range("D19").Value = -33.89637431616226
$ws.Range("E19").Value = 0.0000000000000009596625006765299
$ws.Range("F19").Value = -0.01597402443165519
$ws.Range("G19").Value = -0.01406464665070042
$ws.Range("H19").Value = 0.09147670592520497

$ws.Range("B20").Value = -0.01084563274848796
$ws.Range("C20").Value = 0.0004950637531077616
$ws.Range("D20").Value = -24.59393134766182
$ws.Range("E20").Value = 0.009168122169645156
$ws.Range("F20").Value = -0.01181594302898177
$ws.Range("G20").Value = -0.009875322467994148
$ws.Range("H20").Value = 0.09565040871789482

$ws.Range("B21").Value = -0.008005590703414747
$ws.Range("C21").Value = 0.0004965995510334915
$ws.Range("D21").Value = -18.06034061820105
$ws.Range("E21").Value = 0.002466463359950314
$ws.Range("F21").Value = -0.008978911099275661
$ws.Range("G21").Value = -0.007032270307553834
$ws.Range("H21").Value = 0.09849045076296803

$ws.Range("B22").Value = -0.006442773990630101
$ws.Range("C22").Value = 0.0004922571314374162
$ws.Range("D22").Value = -14.26305225309795
$ws.Range("E22").Value = 0.05551651273030982
$ws.Range("F22").Value = -0.007407583376462374
$ws.Range("G22").Value = -0.005477964604797828
$ws.Range("H22").Value = 0.1000532674757527

$ws.Range("B23").Value = -0.005030046940993602
$ws.Range("C23").Value = 0.0004930780714059177
$ws.Range("D23").Value = -11.23656151001397
$ws.Range("E23").Value = 0.06538108850648955
$ws.Range("F23").Value = -0.005996465341256235
$ws.Range("G23").Value = -0.004063628540730971
$ws.Range("H23").Value = 0.1014659945253892

$ws.Range("B24").Value = -0.004206301800107242
$ws.Range("C24").Value = 0.0004872750687903423
$ws.Range("D24").Value = -9.240635445136027
$ws.Range("E24").Value = 0.03879154343074996
$ws.Range("F24").Value = -0.005161346476548963
$ws.Range("G24").Value = -0.003251257123665521
$ws.Range("H24").Value = 0.1022897396662755

$ws.Range("B25").Value = -0.002419443623081502
$ws.Range("C25").Value = 0.0004773233058582917
$ws.Range("D25").Value = -5.23700088304876
$ws.Range("E25").Value = 0.150637516970459
$ws.Range("F25").Value = -0.003354983142147023
$ws.Range("G25").Value = -0.001483904104015981
$ws.Range("H25").Value = 0.1040765978433013

$ws.Range("B26").Value = 0.0211794654040267
$ws.Range("C26").Value = 0.00114044985415447
$ws.Range("D26").Value = 11.59721952635855
$ws.Range("E26").Value = 0.0114740887632291
$ws.Range("F26").Value = 0.01894421803744472
$ws.Range("G26").Value = 0.02341471277060868
$ws.Range("H26").Value = 0.1276755068704095

